# Commit: "Model returns when valid ID and PIN. Button Resident added to
# transition to Resident Screen"
#
# This TouchGFX "texts.xlsx" workbook has two sheets:
#   Typography  - font/typography definitions (Typography Name, Font, Size, Bpp, ...)
#   Translation - text-id table (TEXT ID, TYPOGRAPHY NAME, ALIGNMENT, DIRECTION, GB)
#
# The edit:
#   1) Renames/repurposes the "Button_Label" typography row into a new,
#      smaller "SmallBtn_Label" typography (font size 25 -> 15) used for
#      the new small "Resident" button.
#   2) Adds new translation rows for the Resident screen: the "Resident
#      Name:" / resident-name labels (sample data "Elon Musk" /
#      "Jeff Bezos") and the new "Resident" button label that transitions
#      to the Resident screen.

$wb = $excel.ActiveWorkbook

$wsTypo  = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# --- Typography sheet ------------------------------------------------
# Row 7 was "Button_Label" / size 25; it becomes "SmallBtn_Label" / size 15.
$wsTypo.Range("B7").Value = "SmallBtn_Label"
$wsTypo.Range("D7").Value = 15

# --- Translation sheet ------------------------------------------------
# Append the new text-id rows (22-25) right after the existing data (row 21).
$newRows = @(
    @{ Row = 22; Id = "SingleUseId21"; Typo = "Default";        Align = "Center"; Dir = "LTR"; Text = "Resident Name:" },
    @{ Row = 23; Id = "SingleUseId22"; Typo = "Default";        Align = "Center"; Dir = "LTR"; Text = "Elon Musk" },
    @{ Row = 24; Id = "SingleUseId24"; Typo = "Default";        Align = "Center"; Dir = "LTR"; Text = "Jeff Bezos" },
    @{ Row = 25; Id = "SingleUseId25"; Typo = "SmallBtn_Label"; Align = "Center"; Dir = "LTR"; Text = "Resident" }
)

foreach ($r in $newRows) {
    $wsTrans.Cells.Item($r.Row, 2).Value = $r.Id    # B - TEXT ID
    $wsTrans.Cells.Item($r.Row, 3).Value = $r.Typo  # C - TYPOGRAPHY NAME
    $wsTrans.Cells.Item($r.Row, 4).Value = $r.Align # D - ALIGNMENT
    $wsTrans.Cells.Item($r.Row, 5).Value = $r.Dir   # E - DIRECTION
    $wsTrans.Cells.Item($r.Row, 6).Value = $r.Text  # F - GB (translated text)
}
